# Workshop fall 2021 update: explicitly stamp every paragraph (body +
# the heading/title styles) with PageBreakBefore = False. The source
# document (a Google Docs export) never wrote <w:pageBreakBefore/> at
# all, so Word's default ("no break") applied implicitly; this edit
# makes that default explicit everywhere, matching the round-tripped
# canonical OOXML produced after a later re-save.

$d = $word.ActiveDocument

# Every paragraph in the document body gets an explicit "no page break
# before" paragraph property.
foreach ($p in $d.Paragraphs) {
    $p.Range.ParagraphFormat.PageBreakBefore = 0
}

# The heading / title paragraph styles (the ones that already carry
# KeepWithNext / KeepLinesTogether) get the same explicit flag added
# to their style definitions.
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", `
                "Heading 5", "Heading 6", "Title", "Subtitle")

foreach ($name in $styleNames) {
    $style = $d.Styles($name)
    $style.ParagraphFormat.PageBreakBefore = 0
}

Write-Output "pageBreakBefore stamped on $($d.Paragraphs.Count) paragraphs and $($styleNames.Count) styles"
